$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "305.22"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.43%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.05"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "6.28%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.009"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.12%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07888"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.02%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.205"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-5.03%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.012"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.68%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.020"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.03%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9230"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.21%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09650"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-4.48%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1891"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.63%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08577"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.69%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03687"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "7.89%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09976"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.73%"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.42%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005638"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.61%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.470"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.11%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "6.97%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3412"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.19%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1318"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.66%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.756"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.69%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2200"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-3.12%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.99%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001233"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.29%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004471"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "3.02%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001400"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "7.81%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "39.95%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01840"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "4.89%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04771"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.50%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.008140"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "4.69%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.82%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007550"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-11.30%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002220"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.95%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01005"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.45%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006282"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "3.19%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.11%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0005802"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.03%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "28.62"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "393.82%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.001721"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-35.95%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00002100"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.11%"
